# LOM3231.xlsx content update
# - shifts/rewrites several label/value cells in rows 10-24
# - adjusts a handful of row heights
# - removes the last row (25), shrinking the sheet from 25 to 24 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height fix-ups -----------------------------------------------
# Rows 17 and 22 must go back to the sheet's default (no explicit custom
# height). Deleting the row and inserting a fresh blank one in its place
# is the simplest way to drop the stored row height.
$ws.Rows(17).Delete()
$ws.Rows(17).Insert()

$ws.Rows(22).Delete()
$ws.Rows(22).Insert()

# Rows that need an explicit custom height.
$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(23).RowHeight = 30

# --- Cell content updates ----------------------------------------------
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("A19").Value = "Critério:"

$ws.Range("A20").Value = "Norma de recuperação:"

$ws.Range("A21").Value = "Bibliografia:"

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
$ws.Range("C23").Value = "LOM3234 -  Óptica Física  (Requisito)`n"

$ws.Range("B24").Value = "LOM3259 -  Materiais e Dispositivos Eletrônicos  (Indicação de Conjunto)`n"
$ws.Range("C24").Value = "LOM3259 -  Materiais e Dispositivos Eletrônicos  (Indicação de Conjunto)`n"

# --- Drop the now-unused trailing row -----------------------------------
$ws.Rows(25).Delete()
